# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
